$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("O1").Value = "age"
$ws.Range("P1").Value = "Zone"
$ws.Range("Q1").Value = "color"

# Data: row, age, zone, color ($null means leave blank but still "touch" the cell)
$data = @(
    @(2,  1577,  "AP2",           "darkolivegreen"),
    @(3,  1593,  "AP2",           "darkolivegreen"),
    @(4,  2375,  $null,           $null),
    @(5,  3596,  "AP1",           "darkolivegreen"),
    @(6,  3694,  "AP1",           "darkolivegreen"),
    @(7,  3792,  "AP1",           "darkolivegreen"),
    @(8,  4150,  "AP1",           "darkolivegreen"),
    @(9,  4361,  "AP1",           "darkolivegreen"),
    @(10, 4675,  $null,           $null),
    @(11, 9224,  "HTM1",          "goldenrod1"),
    @(12, 9337,  "HTM1",          "goldenrod1"),
    @(13, 9398,  "HTM1",          "goldenrod1"),
    @(14, 9467,  "HTM1",          "goldenrod1"),
    @(15, 10933, $null,           $null),
    @(16, 11012, $null,           $null),
    @(17, 11066, $null,           $null),
    @(18, 11184, $null,           $null),
    @(19, 11290, $null,           $null),
    @(20, 11901, "Younger Dryas", "skyblue"),
    @(21, 12187, "Younger Dryas", "skyblue"),
    @(22, 12321, "Younger Dryas", "skyblue"),
    @(23, 12441, "Younger Dryas", "skyblue"),
    @(24, 12630, "Younger Dryas", "skyblue")
)

foreach ($row in $data) {
    $r = $row[0]
    $age = $row[1]
    $zone = $row[2]
    $color = $row[3]

    $ws.Cells.Item($r, 15).Value = $age

    if ($null -ne $zone) {
        $ws.Cells.Item($r, 16).Value = $zone
    } else {
        $ws.Cells.Item($r, 16).Font.Bold = $false
    }

    if ($null -ne $color) {
        $ws.Cells.Item($r, 17).Value = $color
    } else {
        $ws.Cells.Item($r, 17).Font.Bold = $false
    }
}
